# Fruta / hortaliza, semanal
# Inserts two new weekly price rows for "Naranja" (Valencia, Primera) at the
# top of the date-ordered data block, pushing the existing rows 581..670
# down to 583..672 (dimension grows from A1:T670 to A1:T672).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row from 581 downward by two rows, so the two new
# rows land at 581 and 582.
$ws.Rows("581:582").Insert()

# --- New row 581 ---------------------------------------------------------
$ws.Range("A581").Value = 9
$ws.Range("B581").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C581").Value = "Metropolitana"
$ws.Range("D581").Value = 44637
$ws.Range("E581").Value = 13
$ws.Range("F581").Value = "Fruta"
$ws.Range("G581").Value = 100102
$ws.Range("H581").Value = "Cítricos"
$ws.Range("I581").Value = 100102005
$ws.Range("J581").Value = "Naranja"
$ws.Range("K581").Value = "Valencia"
$ws.Range("L581").Value = "Primera"
$ws.Range("M581").Value = 450
$ws.Range("N581").Value = 10000
$ws.Range("O581").Value = 10000
$ws.Range("P581").Value = 10000
$ws.Range("Q581").Value = "`$/caja 18 kilos granel"
$ws.Range("R581").Value = "Región de O'Higgins"
$ws.Range("S581").Value = 556
$ws.Range("T581").Value = 18

# --- New row 582 ---------------------------------------------------------
$ws.Range("A582").Value = 9
$ws.Range("B582").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C582").Value = "Metropolitana"
$ws.Range("D582").Value = 44637
$ws.Range("E582").Value = 13
$ws.Range("F582").Value = "Fruta"
$ws.Range("G582").Value = 100102
$ws.Range("H582").Value = "Cítricos"
$ws.Range("I582").Value = 100102005
$ws.Range("J582").Value = "Naranja"
$ws.Range("K582").Value = "Valencia"
$ws.Range("L582").Value = "Primera"
$ws.Range("M582").Value = 450
$ws.Range("N582").Value = 9000
$ws.Range("O582").Value = 9000
$ws.Range("P582").Value = 9000
$ws.Range("Q582").Value = "`$/malla 18 kilos"
$ws.Range("R582").Value = "Región Metropolitana"
$ws.Range("S582").Value = 500
$ws.Range("T582").Value = 18
